$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E4").Value = "2016-03-24 08:22:50"
$wsZh.Range("H4").Value = "2016-03-24 08:23:19"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E4").Value = "2016-03-24 08:22:55"
$wsDe.Range("H4").Value = "2016-03-24 08:23:26"
